# Fruta / hortaliza, semanal
# Insert 5 new weekly price rows for "Murcott" mandarins (Provincia de Quillota /
# Provincia de San Felipe de Aconcagua) right before the existing row 275,
# pushing the existing data (previously rows 275:329) down to rows 280:334.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 5 blank rows above the current row 275.
$ws.Range("A275:T279").Insert()

# --- Row 275 ---
$ws.Range("A275").Value2 = 9
$ws.Range("B275").Value2 = "Vega Central Mapocho de Santiago"
$ws.Range("C275").Value2 = "Metropolitana"
$ws.Range("D275").Value2 = 44476
$ws.Range("E275").Value2 = 13
$ws.Range("F275").Value2 = "Fruta"
$ws.Range("G275").Value2 = 100102
$ws.Range("H275").Value2 = "Cítricos"
$ws.Range("I275").Value2 = 100102004
$ws.Range("J275").Value2 = "Mandarina"
$ws.Range("K275").Value2 = "Murcott"
$ws.Range("L275").Value2 = "Especial"
$ws.Range("M275").Value2 = 300
$ws.Range("N275").Value2 = 5000
$ws.Range("O275").Value2 = 5000
$ws.Range("P275").Value2 = 5000
$ws.Range("Q275").Value2 = "`$/bandeja 10 kilos"
$ws.Range("R275").Value2 = "Provincia de Quillota"
$ws.Range("S275").Value2 = 500
$ws.Range("T275").Value2 = 10

# --- Row 276 ---
$ws.Range("A276").Value2 = 9
$ws.Range("B276").Value2 = "Vega Central Mapocho de Santiago"
$ws.Range("C276").Value2 = "Metropolitana"
$ws.Range("D276").Value2 = 44476
$ws.Range("E276").Value2 = 13
$ws.Range("F276").Value2 = "Fruta"
$ws.Range("G276").Value2 = 100102
$ws.Range("H276").Value2 = "Cítricos"
$ws.Range("I276").Value2 = 100102004
$ws.Range("J276").Value2 = "Mandarina"
$ws.Range("K276").Value2 = "Murcott"
$ws.Range("L276").Value2 = "Extra (doble especial)"
$ws.Range("M276").Value2 = 250
$ws.Range("N276").Value2 = 6000
$ws.Range("O276").Value2 = 6000
$ws.Range("P276").Value2 = 6000
$ws.Range("Q276").Value2 = "`$/bandeja 10 kilos"
$ws.Range("R276").Value2 = "Provincia de Quillota"
$ws.Range("S276").Value2 = 600
$ws.Range("T276").Value2 = 10

# --- Row 277 ---
$ws.Range("A277").Value2 = 9
$ws.Range("B277").Value2 = "Vega Central Mapocho de Santiago"
$ws.Range("C277").Value2 = "Metropolitana"
$ws.Range("D277").Value2 = 44476
$ws.Range("E277").Value2 = 13
$ws.Range("F277").Value2 = "Fruta"
$ws.Range("G277").Value2 = 100102
$ws.Range("H277").Value2 = "Cítricos"
$ws.Range("I277").Value2 = 100102004
$ws.Range("J277").Value2 = "Mandarina"
$ws.Range("K277").Value2 = "Murcott"
$ws.Range("L277").Value2 = "Primera"
$ws.Range("M277").Value2 = 200
$ws.Range("N277").Value2 = 4000
$ws.Range("O277").Value2 = 4000
$ws.Range("P277").Value2 = 4000
$ws.Range("Q277").Value2 = "`$/bandeja 10 kilos"
$ws.Range("R277").Value2 = "Provincia de Quillota"
$ws.Range("S277").Value2 = 400
$ws.Range("T277").Value2 = 10

# --- Row 278 ---
$ws.Range("A278").Value2 = 9
$ws.Range("B278").Value2 = "Vega Central Mapocho de Santiago"
$ws.Range("C278").Value2 = "Metropolitana"
$ws.Range("D278").Value2 = 44476
$ws.Range("E278").Value2 = 13
$ws.Range("F278").Value2 = "Fruta"
$ws.Range("G278").Value2 = 100102
$ws.Range("H278").Value2 = "Cítricos"
$ws.Range("I278").Value2 = 100102004
$ws.Range("J278").Value2 = "Mandarina"
$ws.Range("K278").Value2 = "Murcott"
$ws.Range("L278").Value2 = "Primera"
$ws.Range("M278").Value2 = 22
$ws.Range("N278").Value2 = 130000
$ws.Range("O278").Value2 = 140000
$ws.Range("P278").Value2 = 135455
$ws.Range("Q278").Value2 = "`$/bins (400 kilos)"
$ws.Range("R278").Value2 = "Provincia de San Felipe de Aconcagua"
$ws.Range("S278").Value2 = 339
$ws.Range("T278").Value2 = 400

# --- Row 279 ---
$ws.Range("A279").Value2 = 9
$ws.Range("B279").Value2 = "Vega Central Mapocho de Santiago"
$ws.Range("C279").Value2 = "Metropolitana"
$ws.Range("D279").Value2 = 44476
$ws.Range("E279").Value2 = 13
$ws.Range("F279").Value2 = "Fruta"
$ws.Range("G279").Value2 = 100102
$ws.Range("H279").Value2 = "Cítricos"
$ws.Range("I279").Value2 = 100102004
$ws.Range("J279").Value2 = "Mandarina"
$ws.Range("K279").Value2 = "Murcott"
$ws.Range("L279").Value2 = "Segunda"
$ws.Range("M279").Value2 = 250
$ws.Range("N279").Value2 = 3000
$ws.Range("O279").Value2 = 3000
$ws.Range("P279").Value2 = 3000
$ws.Range("Q279").Value2 = "`$/bandeja 10 kilos"
$ws.Range("R279").Value2 = "Provincia de Quillota"
$ws.Range("S279").Value2 = 300
$ws.Range("T279").Value2 = 10
